$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E5").Value = "fullRNASEQ"

$ws.Range("D6:F15").Select()
